$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Price" column (D) holds numbers stored as literal text in the
# source data. Values containing two '.' separators (e.g. "97.655.27") are
# never misread as numbers, but plain decimal-looking values (e.g. "623.43")
# must be forced to stay text (matching the original inlineStr cell type)
# by switching the cell to the text number-format right before assigning.

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "97.655.27"
$ws.Range("E2").Value = "  +3.66%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.337.12"

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 (Solana)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.56"
$ws.Range("E5").Value = "  +3.66%  "

# Row 6 (BNB)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.43"
$ws.Range("E6").Value = "  +1.41%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  +0.40%  "

# Row 8 (Dogecoin)
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 (USDC)
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 (LidoStakedEther)
$ws.Range("D10").Value = "3.335.70"
$ws.Range("E10").Value = "  +7.58%  "

# Row 11 (Cardano)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.800"
$ws.Range("E11").Value = "  -3.47%  "

# Row 12 (TRON)
$ws.Range("E12").Value = "  +1.45%  "

# Row 13 (WrappedBTC)
$ws.Range("D13").Value = "97.390.05"
$ws.Range("E13").Value = "  +3.92%  "

# Row 14 & 15 swap (Avalanche <-> ShibaInu)
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.72"
$ws.Range("E15").Value = "  +2.48%  "

# Row 16 (WrappedliquidstakedEther2.0)
$ws.Range("D16").Value = "3.956.47"
$ws.Range("E16").Value = "  +7.45%  "

# Row 17 (Toncoin)
$ws.Range("E17").Value = "  +2.32%  "

# Row 18 (WrappedEther)
$ws.Range("D18").Value = "3.331.77"
$ws.Range("E18").Value = "  +7.52%  "

# Row 19 (SuiNetwork)
$ws.Range("E19").Value = "  +0.23%  "

# Row 20 (Chainlink)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  +2.76%  "

# Row 21 (BitcoinCash)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "495.33"
$ws.Range("E21").Value = "  +11.22%  "

# Row 22 & 23 swap (Polkadot <-> PEPE)
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000213"
$ws.Range("E22").Value = "  +6.42%  "

$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").Value = "  -0.84%  "

# Row 24 (Uniswap)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.27"
$ws.Range("E24").Value = "  +3.22%  "

# Row 25 (NEARProtocol)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26 (Litecoin)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.82"
$ws.Range("E26").Value = "  +3.20%  "

# Row 27 (Aptos)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("E27").Value = "  -0.68%  "

# Row 28 (WrappedeETH)
$ws.Range("D28").Value = "3.504.11"
$ws.Range("E28").Value = "  +7.18%  "

# Row 29 (Dai)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30 (Cronos)
$ws.Range("E30").Value = "  +1.33%  "

# Row 31 (Stellar)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.242"
$ws.Range("E31").Value = "  -1.27%  "

# Row 32 (Hedera)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -1.13%  "

# Row 33 (Binance-PegBSC-USD)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 (InternetComputer(DFINITY))
$ws.Range("E34").Value = "  +1.14%  "

# Row 35 (EthereumClassic)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.71"
$ws.Range("E35").Value = "  +6.47%  "

# Row 36 (Kaspa)
$ws.Range("E36").Value = "  -4.91%  "

# Row 37 (RenderToken)
$ws.Range("E37").Value = "  -4.32%  "

# Row 38 (Bittensor)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "504.88"
$ws.Range("E38").Value = "  +5.95%  "

# Row 39 (PancakeSwap)
$ws.Range("E39").Value = "  +2.29%  "

# Row 40 (WhiteBITCoin)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.73"
$ws.Range("E40").Value = "  +3.14%  "

# Row 41 (PolygonEcosystemToken)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.451"
$ws.Range("E41").Value = "  -0.06%  "

# Row 42 (Fetch.AI)
$ws.Range("E42").Value = "  +0.27%  "

# Row 43 (dogwifhat)
$ws.Range("E43").Value = "  +1.61%  "

# Row 44, 45, 46 rotate: MantraDAO/USDe/ARBITRUM -> ARBITRUM/MantraDAO/USDe
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.792"
$ws.Range("E44").Value = "  +15.30%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  -7.78%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47 (Monero)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.15"
$ws.Range("E47").Value = "  +0.20%  "

# Row 48 (Stacks)
$ws.Range("E48").Value = "  +5.68%  "

# Row 49 (Filecoin)
$ws.Range("E49").Value = "  +3.96%  "

# Row 50 (VeChain)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +3.80%  "

# Row 51 (OKB)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.92"
$ws.Range("E51").Value = "  +2.54%  "
